$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new transaction rows (24-26) below the existing data,
# following the same E/N/P/T layout used by the rows above.
$newRows = @(
    @("Deposit",    "Crypto", "ETH", 628.65880000000004),
    @("Withdrawal", "Crypto", "ETH", 1000),
    @("Deposit",    "Crypto", "ETH", 5356.6178)
)

$r = 24
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 5).Value  = $row[0]   # column E - Transaction Type
    $ws.Cells.Item($r, 14).Value = $row[1]   # column N - Payment Type
    $ws.Cells.Item($r, 16).Value = $row[2]   # column P - InternalComment
    $ws.Cells.Item($r, 20).Value = $row[3]   # column T - USD Amount
    $r++
}

# Reflect the user's scroll position / selection at the time of saving.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F20:M31").Select()
